# 2023.01.21  basically finished for: editting support for PredefTrainFilter
# Add three new rows to the "操作命令表" (command list) sheet describing the
# new PredefTrainFilter-related commands (Add/Remove/Update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateFormat = "yyyy""年""m""月""d""日"";@"

$rows = @(
    @{ No = 57; Name = "新建预置列车筛选器"; Cmd = "AddTrainFilter";    Support = "支持"; Merge = "否"; Date = 44947 },
    @{ No = 58; Name = "删除预置列车筛选器"; Cmd = "RemoveTrainFilter"; Support = "支持"; Merge = "否"; Date = 44947 },
    @{ No = 59; Name = "更新列车筛选器";     Cmd = "UpdateTrainFilter"; Support = "支持"; Merge = "否"; Date = 44947 }
)

$r = 59
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.No
    $ws.Cells.Item($r, 2).Value = $row.Name
    $ws.Cells.Item($r, 3).Value = $row.Cmd
    $ws.Cells.Item($r, 4).Value = $row.Support
    $ws.Cells.Item($r, 5).Value = $row.Merge
    $ws.Cells.Item($r, 6).Value = $row.Date
    $ws.Cells.Item($r, 6).NumberFormat = $dateFormat
    $r++
}

$ws.Range("G61").Select() | Out-Null
